# Timesheet.xlsx - "Added sum on timesheet"
#
# Adds a small "person hours" summary block in column J:
#   J3 - label  "Estimate of person hours for project 3:"
#   J4 - value  120 (the estimate), yellow-filled, centered
#   J7 - label  "Actual Person hours for project 3:"
#   J8 - array formula totalling both weeks' hours, yellow-filled
# Also fills in the previously-blank H13 (Freeman/Sunday, week 2) with 2
# hours, widens column J so the labels are readable, and leaves the
# selection where the author's cursor ended up (J16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- missing data point on the second week's timesheet grid ---
$ws.Range("H13").Value = 2

# --- estimate label + value (J3:J4) ---
$ws.Range("J4").Value = 120
$ws.Range("J4").Interior.Color = 65535      # yellow
$ws.Range("J4").HorizontalAlignment = -4108 # xlCenter
$ws.Range("J4").VerticalAlignment = -4108   # xlCenter

$ws.Range("J3").Value = "Estimate of person hours for project 3:"
$ws.Range("J3").Font.Underline = 2          # xlUnderlineStyleSingle
$ws.Range("J3").HorizontalAlignment = -4108 # xlCenter
$ws.Range("J3").VerticalAlignment = -4108   # xlCenter

# --- actual label + totalling array formula (J7:J8) ---
$ws.Range("J7").Value = "Actual Person hours for project 3:"
$ws.Range("J7").Font.Underline = 2          # xlUnderlineStyleSingle

$ws.Range("J8").FormulaArray = "=SUM(B3:H7+B9:H13)"
$ws.Range("J8").Interior.Color = 65535      # yellow

# --- column sizing + final selection, matching the saved workbook state ---
$ws.Columns("J").ColumnWidth = 29.5

$ws.Range("J16").Select() | Out-Null
